# Daily auto push: append a new scraped data point for 2026/01/25 onto the
# existing run of 2026/01/25 rows (row 696), shifting every following row
# (old 696..737, the 2026/12/29 .. 2027/01/05 block) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 696 (pushes old 696..737 down to 697..738)
$ws.Rows.Item(696).Insert()

# Fill the new row with the latest scraped data point.
# Column A stores dates as plain text (not Excel date serials) in this sheet,
# so force text entry with a leading apostrophe and strip the resulting
# "quote prefix" formatting so the cell style matches its plain-text siblings.
$ws.Range("A696").Value = "'2026/01/25"
$ws.Range("A696").ClearFormats()
$ws.Range("B696").Value = "日"
$ws.Range("C696").Value = 19
$ws.Range("D696").Value = 17
